$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Air_Sabre_RoundTrip")

# Update the AirLines column (L2:L5) value from "Southwest Airlines" to "SOUTHWEST"
$ws.Range("L2:L5").Value = "SOUTHWEST"

# Auto-fit column L to the new content (matches bestFit width observed in target)
$ws.Columns.Item(12).EntireColumn.AutoFit()

# Reproduce the final selection / viewport state
$ws.Range("L3:L5").Select()
$excel.ActiveWindow.ScrollColumn = 3
